# Generate Report for Handoff
# Updates the localization-status workbook: flips the two tracked files from
# "handed back" to "ready for handoff", refreshes their timestamps, switches
# the zh-cn priority to machine translation, and records a handback-version
# warning for the cd80da35 file. Also mirrors the cosmetic column-width
# tweaks that Excel applied when the report was regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "Ready for handoff"

$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# ---------------------------------------------------------------------------
# Overview "Latest HO Xliff Generate Date" column (shared with de-de's
# "Latest Handoff Datetime" column - both pointed at the same text before).
$wsOverview.Range("G2").Value = "2016-09-07 10:15:47"
$wsOverview.Range("G3").Value = "2016-09-07 10:15:47"

$wsDe.Range("H2").Value = "2016-09-07 10:15:47"
$wsDe.Range("H3").Value = "2016-09-07 10:15:47"

# zh-cn "Latest Handoff Datetime"
$wsZh.Range("H2").Value = "2016-09-07 10:15:35"
$wsZh.Range("H3").Value = "2016-09-07 10:15:35"

# ---------------------------------------------------------------------------
# Priority: "ht" (human translation) -> "mt" (machine translation)
# ---------------------------------------------------------------------------
$wsZh.Range("E2").Value = "mt"
$wsZh.Range("E3").Value = "mt"

$wsDe.Range("E2").Value = "mt"
$wsDe.Range("E3").Value = "mt"

# ---------------------------------------------------------------------------
# New "Error Detail" message on the zh-cn sheet for the cd80da35 handback
# ---------------------------------------------------------------------------
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7258dcbb96e03d0fb38bda1d58fcb1d9d9c34908/e2e/cd80da35-0eab-46fd-9ac2-4d0602d92db7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48bfdf6eb946816806e0ec189cbd1b28b3b1031b/e2e/cd80da35-0eab-46fd-9ac2-4d0602d92db7.md."

# ---------------------------------------------------------------------------
# Column width tweaks (cosmetic re-layout that accompanied the report
# regeneration). The COM ColumnWidth setter here snaps to whole-pixel
# increments, so we feed it the pixel-equivalent input that lands closest to
# the target stored width.
# ---------------------------------------------------------------------------
$padding = 5.0 / 6.0

$wsOverview.Columns.Item(5).ColumnWidth = (17.2159881591797 - $padding)
$wsOverview.Columns.Item(6).ColumnWidth = (17.2159881591797 - $padding)

$wsZh.Columns.Item(3).ColumnWidth = (17.2159881591797 - $padding)
$wsZh.Columns.Item(16).ColumnWidth = (40 - $padding)

$wsDe.Columns.Item(3).ColumnWidth = (17.2159881591797 - $padding)
$wsDe.Columns.Item(16).ColumnWidth = (40 - $padding)
